$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old per-customer product columns (B:F) for rows 1-4.
# Column A (customer emails) and its hyperlinks/styles stay untouched.
$ws.Range("B1:F4").ClearContents()

# Row 2 values: append " cases" to the existing G2/H2/I2 order text.
$ws.Range("G2").Value = "1 P2FETT cases"
$ws.Range("H2").Value = "1 P1PAPP cases"
$ws.Range("I2").Value = "5 P2LING cases"

# Row 1/2 new "lbs" column.
$ws.Range("J1").Value = "2 lbs Egg Fettuccine"
$ws.Range("J2").Value = "2 P1FETT lbs"

# Update the active selection to match the saved workbook state.
$ws.Range("K4").Select()
